$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows (bottom-to-top using original row numbers so positions stay valid)
$ws.Rows.Item(24).Insert()   # will hold RECLTD after second insert below
$ws.Rows.Item(24).Insert()   # will hold ONGC
$ws.Rows.Item(10).Insert()   # will hold DIXON
$ws.Rows.Item(5).Insert()    # will hold ASHOKLEY
$ws.Rows.Item(3).Insert()    # will hold ABB

# Populate every data row with final ticker/value layout
$ws.Range("A2").Value = "AARTIIND"
$ws.Range("B2").Value = 628.7
$ws.Range("C2").Value = 640.75

$ws.Range("A3").Value = "ABB"
$ws.Range("B3").Value = 7876.5
$ws.Range("C3").ClearContents()

$ws.Range("A4").Value = "ADANI"
$ws.Range("B4").NumberFormat = "#,##0.00"
$ws.Range("B4").Value = 3020.15
$ws.Range("C4").Value = 3028

$ws.Range("A5").Value = "APOLLO"
$ws.Range("B5").Value = 492.35
$ws.Range("C5").Value = 497.5

$ws.Range("A6").Value = "ASHOKLEY"
$ws.Range("B6").Value = 253.8
$ws.Range("C6").ClearContents()

$ws.Range("A7").Value = "BAJFINSV"
$ws.Range("B7").NumberFormat = "#,##0.00"
$ws.Range("B7").Value = 1755.65
$ws.Range("C7").Value = 1713.5

$ws.Range("A8").Value = "BAJFIN"
$ws.Range("B8").NumberFormat = "#,##0.00"
$ws.Range("B8").Value = 7063.55
$ws.Range("C8").Value = 6900

$ws.Range("A9").Value = "BANBK"
$ws.Range("B9").Value = 193.93
$ws.Range("C9").Value = 196.66

$ws.Range("A10").Value = "BARODA"
$ws.Range("B10").Value = 249.8
$ws.Range("C10").Value = 249.85

$ws.Range("A11").Value = "BN"
$ws.Range("B11").NumberFormat = "#,##0.00"
$ws.Range("B11").Value = 51490.75
$ws.Range("C11").Value = 51143.85

$ws.Range("A12").Value = "DIXON"
$ws.Range("B12").NumberFormat = "#,##0.00"
$ws.Range("B12").Value = 13201.9
$ws.Range("C12").ClearContents()

$ws.Range("A13").Value = "DLF"
$ws.Range("B13").Value = 831.9
$ws.Range("C13").Value = 837.1

$ws.Range("A14").Value = "EICHER"
$ws.Range("B14").NumberFormat = "#,##0.00"
$ws.Range("B14").Value = 4915.9
$ws.Range("C14").Value = 4946.85

$ws.Range("A15").Value = "ESCORTS"
$ws.Range("B15").NumberFormat = "#,##0.00"
$ws.Range("B15").Value = 3816.5
$ws.Range("C15").Value = 3854.75

$ws.Range("A16").Value = "FEDBANK"
$ws.Range("B16").Value = 195.92
$ws.Range("C16").Value = 195.56

$ws.Range("A17").Value = "HCL"
$ws.Range("B17").NumberFormat = "#,##0.00"
$ws.Range("B17").Value = 1751.85
$ws.Range("C17").Value = 1719.45

$ws.Range("A18").Value = "HINDALCO"
$ws.Range("B18").Value = 700.5
$ws.Range("C18").Value = 705.05

$ws.Range("A19").Value = "IGL"
$ws.Range("B19").Value = 543.35
$ws.Range("C19").Value = 538.7

$ws.Range("A20").Value = "INDUSIND"
$ws.Range("B20").NumberFormat = "#,##0.00"
$ws.Range("B20").Value = 1417.45
$ws.Range("C20").Value = 1415.75

$ws.Range("A21").Value = "JIND"
$ws.Range("B21").Value = 960.5
$ws.Range("C21").Value = 968.9

$ws.Range("A22").Value = "LIC"
$ws.Range("B22").Value = 675.95
$ws.Range("C22").Value = 674.9

$ws.Range("A23").Value = "M&M"
$ws.Range("B23").NumberFormat = "#,##0.00"
$ws.Range("B23").Value = 2757.6
$ws.Range("C23").Value = 2798

$ws.Range("A24").Value = "M&MFIN"
$ws.Range("B24").Value = 313.4
$ws.Range("C24").Value = 315.25

$ws.Range("A25").Value = "NIFTY"
$ws.Range("B25").NumberFormat = "#,##0.00"
$ws.Range("B25").Value = 25265.2
$ws.Range("C25").Value = 25048.35

$ws.Range("A26").Value = "NTPC"
$ws.Range("B26").Value = 409.9
$ws.Range("C26").Value = 409.05

$ws.Range("A27").Value = "ONGC"
$ws.Range("B27").Value = 329.6
$ws.Range("C27").ClearContents()

$ws.Range("A28").Value = "RECLTD"
$ws.Range("B28").Value = 626.25
$ws.Range("C28").ClearContents()

$ws.Range("A29").Value = "SBIN"
$ws.Range("B29").Value = 814.5
$ws.Range("C29").Value = 809.4

$ws.Range("A30").Value = "SUNTV"
$ws.Range("B30").Value = 809.9
$ws.Range("C30").Value = 815.25

$ws.Range("A31").Value = "TM"
$ws.Range("B31").NumberFormat = "#,##0.00"
$ws.Range("B31").Value = 1121.65
$ws.Range("C31").Value = 1074.55

$ws.Range("A32").Value = "TP"
$ws.Range("B32").Value = 430.9
$ws.Range("C32").Value = 431.4

$ws.Range("A33").Value = "TS"
$ws.Range("B33").Value = 152.97
$ws.Range("C33").Value = 153.7

$ws.Range("A34").Value = "VEDL"
$ws.Range("B34").Value = 463.4
$ws.Range("C34").Value = 466.05

# Restore selection to match target state
$ws.Range("H12").Select()
